# Clean up and Page creator
#
# On the "ListeHTML" sheet: the "<p>" tag row ("ecris, paragraphe") is
# relabelled to "paragraphe" and a second, duplicate "paragraphe" / "<p>"
# row is inserted right after it (rows 38-52 shift down to 39-53).
# The saved scroll position / selection are also updated to match where the
# author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ListeHTML")
$ws.Activate()

# Insert a new row at position 38 (shifts old rows 38-52 down to 39-53)
$ws.Rows.Item(38).Insert()

# Row 37: rename label from "ecris, paragraphe" to "paragraphe" (tag stays <p>)
$ws.Range("A37").Value = "paragraphe"
$ws.Range("B37").Value = "<p>"

# New row 38: duplicate "paragraphe" / "<p>" entry
$ws.Range("A38").Value = "paragraphe"
$ws.Range("B38").Value = "<p>"

# Row heights for the two "paragraphe" rows
$ws.Rows.Item(37).RowHeight = 16.2
$ws.Rows.Item(38).RowHeight = 15.6

# Restore the view: scrolled position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("C38").Select()
